$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B11: rule name changes from "R40" to "1"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
